# "Update eggs and weights" -- append three new days of logged data
# (2021-02-26, 2021-02-27, 2021-02-28) to the bottom of the tracking sheet,
# then leave the view scrolled/selected where the user last left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the last existing data row (A38) down onto
# the three new rows so the new dates pick up the same number format/style
# instead of minting a brand-new style entry.
$ws.Range("A38").Copy()
$ws.Range("A39:A41").PasteSpecial(-4122)

# Row 39 -- 2021-02-26 (all zero / rest day)
$ws.Range("A39").Value = 44253
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0

# Row 40 -- 2021-02-27
$ws.Range("A40").Value = 44254
$ws.Range("B40").Value = 60
$ws.Range("C40").Value = 45
$ws.Range("D40").Value = 12
$ws.Range("E40").Value = 30
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0

# Row 41 -- 2021-02-28
$ws.Range("A41").Value = 44255
$ws.Range("B41").Value = 24
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0

# Match the saved view state: selection on J40, scrolled down so row 25 is
# at the top of the visible pane.
$ws.Range("J40").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
